$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Updated test data values for Battery Standby / Battery Alarm Limits
$ws.Range("N5").Value = 0.285
$ws.Range("O5").Value = 0.444
$ws.Range("N6").Value = 0.285
$ws.Range("O6").Value = 0.444

# Make "Add Panels" the active sheet and set the view/selection state
$ws.Activate()
$ws.Range("N6:O6").Select()
